# This script refreshes the "cryptos" price table (columns B:E, rows 2-51)
# with the latest scraped values, matching the commit
# "Updated cryptos list ... with GitHub Actions".
#
# All price/volume cells in the sheet are stored as literal text (the
# original file uses <c t="inlineStr"> cells, not numbers). Values that
# look numeric (e.g. "114.49", "1.00", "173.70") are written with a
# leading apostrophe so Excel keeps them as text instead of silently
# re-typing them as numbers - which would, for example, drop the trailing
# zero in "173.70" -> 173.7, or turn "1.00" into 1. ClearFormats()
# afterwards removes the residual "quote prefix" cell style that trick
# leaves behind, so the cell keeps its original (General/default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $range = $ws.Range($addr)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $range.Value = "'" + $text
        $range.ClearFormats()
    } else {
        $range.Value = $text
    }
}

Set-TextValue "D2" '43.697.56'
Set-TextValue "E2" '  -0.04%  '

Set-TextValue "D3" '2.284.94'

Set-TextValue "E4" '  -0.47%  '

Set-TextValue "D5" '114.49'
Set-TextValue "E5" '  +1.68%  '

Set-TextValue "D6" '265.39'
Set-TextValue "E6" '  -1.66%  '

Set-TextValue "E7" '  +2.67%  '

Set-TextValue "E8" '  +0.09%  '

Set-TextValue "D9" '0.614'
Set-TextValue "E9" '  -1.06%  '

Set-TextValue "D10" '47.24'
Set-TextValue "E10" '  -1.57%  '

Set-TextValue "D11" '0.0938'
Set-TextValue "E11" '  -1.21%  '

Set-TextValue "E12" '  +0.83%  '

Set-TextValue "E13" '  +1.47%  '

Set-TextValue "D14" '15.42'
Set-TextValue "E14" '  -2.46%  '

Set-TextValue "D15" '2.631.58'
Set-TextValue "E15" '  -0.30%  '

Set-TextValue "D16" '0.873'
Set-TextValue "E16" '  +2.59%  '

Set-TextValue "D17" '2.286.64'
Set-TextValue "E17" '  -0.25%  '

Set-TextValue "D18" '43.549.73'
Set-TextValue "E18" '  -0.42%  '

Set-TextValue "E19" '  +0.43%  '

Set-TextValue "D20" '6.83'
Set-TextValue "E20" '  +1.68%  '

Set-TextValue "D21" '72.42'
Set-TextValue "E21" '  +0.14%  '

Set-TextValue "D22" '2.44'
Set-TextValue "E22" '  -0.31%  '

Set-TextValue "D23" '236.35'
Set-TextValue "E23" '  +1.52%  '

Set-TextValue "D24" '9.47'
Set-TextValue "E24" '  -3.66%  '

Set-TextValue "E25" '  +0.97%  '

Set-TextValue "E26" '  +1.80%  '

Set-TextValue "E27" '  -1.19%  '

Set-TextValue "D28" '41.81'
Set-TextValue "E28" '  -0.25%  '

Set-TextValue "D29" '3.40'
Set-TextValue "E29" '  +0.09%  '

Set-TextValue "E30" '  -0.99%  '

Set-TextValue "D31" '173.70'
Set-TextValue "E31" '  -1.06%  '

Set-TextValue "D32" '21.70'
Set-TextValue "E32" '  +0.75%  '

Set-TextValue "D33" '0.0910'
Set-TextValue "E33" '  -2.06%  '

Set-TextValue "D34" '5.70'

Set-TextValue "E35" '  +1.65%  '

Set-TextValue "E36" '  +5.76%  '

Set-TextValue "D37" '4.68'
Set-TextValue "E37" '  +0.28%  '

Set-TextValue "D38" '3.94'
Set-TextValue "E38" '  +3.01%  '

Set-TextValue "D39" '0.105'
Set-TextValue "E39" '  -1.96%  '

Set-TextValue "E40" '  +6.74%  '

Set-TextValue "D41" '14.34'
Set-TextValue "E41" '  +4.24%  '

Set-TextValue "D42" '74.24'
Set-TextValue "E42" '  +0.01%  '

Set-TextValue "E43" '  -3.12%  '

Set-TextValue "B44" 'FirstDigitalUSD'
Set-TextValue "C44" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D44" '1.00'
Set-TextValue "E44" '  -0.12%  '

Set-TextValue "B45" 'THORChain'
Set-TextValue "C45" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D45" '5.96'
Set-TextValue "E45" '  -6.07%  '

Set-TextValue "E46" '  -1.35%  '

Set-TextValue "E47" '  +4.05%  '

Set-TextValue "D48" '73.78'
Set-TextValue "E48" '  +34.74%  '

Set-TextValue "E49" '  -2.49%  '

Set-TextValue "E50" '  +0.18%  '

Set-TextValue "D51" '100.44'
Set-TextValue "E51" '  -2.81%  '
